# Added new test to Authoring suite.
#
# Appends a new test-case row (row 40) to the bottom of the "Test Cases"
# sheet (the active sheet, tabSelected) describing the new
# AppreciateUnAppreciateOthersPost test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 already carries the exact formatting pattern ("Y"-style data row
# without the wrap-text wrapper on column D) this new row should have, so
# copy its formats down onto the new row before filling in the values.
$ws.Range("A31:E31").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A40").Value = "AppreciateUnAppreciateOthersPost"
$ws.Range("B40").Value = "OPQA-342|OPQA-359"
$ws.Range("C40").Value = "Verify that user is able to Appreciate/Un Appreciate their others post"
$ws.Range("D40").Value = "Y"

# Scroll / reselect onto the freshly-entered row, same as Excel leaves the
# view after typing a new last row of data.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A40").Select()
